$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Frequency of data input - same value for all entries" table
#    cell currently holds the text split across three runs; collapse
#    it down to a single run/text by replacing the phrase with itself
#    - Word's Find/Replace rewrites the matched range as one run.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Frequency of data input - same value for all entries",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Frequency of data input - same value for all entries", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Move the blank separator paragraph that currently sits right
#    after "...countries in the European Community." (there are two
#    blank paragraphs there in a row) so it instead sits right before
#    the "Research Question" Heading3 paragraph.
# ------------------------------------------------------------------

# Locate the "Research Question" heading paragraph.
$researchQ = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Research Question") {
        $researchQ = $p
        break
    }
}

# InsertParagraphBefore() splits the paragraph's range and inserts a
# new blank paragraph ahead of it. The $researchQ object reference
# rebinds to that new (blank) paragraph afterwards - since it's now
# the paragraph occupying the original slot - so style it directly.
$researchQ.Range.InsertParagraphBefore() | Out-Null
$researchQ.Style = "Normal"
$researchQ.Format.LeftIndent = 21.3

# Now remove one of the two blank paragraphs following "...European
# Community." (the first of the pair), leaving the second blank
# paragraph in place right before the "Data Understanding" heading.
$europeanCommunity = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "countries in the European Community\.") {
        $europeanCommunity = $p
        break
    }
}

$blankToRemove = $europeanCommunity.Next()
$blankToRemove.Range.Delete() | Out-Null

Write-Host "Edit complete"
